$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.179.55'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.420.85'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.00'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.10'
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.91'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").Value = '2.852.25'
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").Value = '60.081.04'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("D17").Value = '2.499.14'
$ws.Range("E17").Value = '  +3.20%  '
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.49'
$ws.Range("E19").Value = '  +2.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '328.00'
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.39'
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("E24").Value = '  +2.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.68'
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("D28").Value = '0.0₃0773'
$ws.Range("E28").Value = '  -2.12%  '
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.22'
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.10'
$ws.Range("E31").Value = '  -3.37%  '
$ws.Range("E32").Value = '  +1.42%  '
$ws.Range("E33").Value = '  -4.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.54'
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("E36").Value = '  +1.07%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.21'
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '329.03'
$ws.Range("E39").Value = '  +2.03%  '
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '145.44'
$ws.Range("E41").Value = '  +4.13%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.66'
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.03'
$ws.Range("E43").Value = '  +2.43%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0967'
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0515'
$ws.Range("E45").Value = '  -1.31%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.574'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0224'
$ws.Range("E47").Value = '  -1.39%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.04'
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.58'
$ws.Range("E49").Value = '  -2.83%  '
$ws.Range("B50").Value = 'ZEEBU'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.65'
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("B51").Value = 'BitgetToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.946'
$ws.Range("E51").Value = '  -0.51%  '
